$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up the "District" column (G) for these rows: it previously held
# the school name concatenated with "Gadag"; it should just say "Gadag".
$ws.Range("G15").Value = "Gadag"
$ws.Range("G28").Value = "Gadag"
$ws.Range("G41").Value = "Gadag"
$ws.Range("G48").Value = "Gadag"
